$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) with the same layout/style as row 3.
# Copy the formatting (cell style / number format) from row 3 first so the
# new date & boolean cells reuse the existing style records instead of
# minting new ones.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(4, 1).Value = 42633.676759259259
$ws.Cells.Item(4, 2).Value = $false
$ws.Cells.Item(4, 3).Value = 9948
$ws.Cells.Item(4, 4).Value = 10000
$ws.Cells.Item(4, 5).Value = 19.32
$ws.Cells.Item(4, 6).Value = 19.12
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = -1.04
$ws.Cells.Item(4, 9).Value = $false
